$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column H
$ws.Range("H1").Value = "NRMSE"

# Updated data values (rows 2-6, columns A-H)
$data = @(
    @(4.4210263244707022, 3.1546950902821074, 0.41146286843744595, 0.64145371496113879, 0.54316375521386151, 0.58853713156255405, 0.77377551794237132, 0.28087841959788457),
    @(4.9643016717198858, 4.1375026227856795, 1.1106234942423863, 1.0538612310178159, 1.107468582116081, -0.11062349424238627, 0.0022284809797388123, 0.37102404123467009),
    @(2.9115312232372008, 2.6564621777071999, 0.50226418826640706, 0.7087059956472832, 0.76687707208637412, 0.49773581173359294, 0.7615691800612383, 0.24303265636370625),
    @(3.6717651553544939, 3.3499041146440005, 0.85147151349506167, 0.92275214087807012, 1.0654911306119592, 0.14852848650493833, 0.94024265242527749, 0.25012024219036066),
    @(2.1198525763532605, 1.4446834055963067, 0.14116097509208772, 0.37571395381604838, 0.31460875557410856, 0.85883902490791231, 0.93529252501144367, 0.12573265577421472)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}

# Column width adjustments (closest achievable values given COM rounding to 1/6 character units)
$ws.Columns.Item(6).ColumnWidth = 12.666666666666666
$ws.Columns.Item(7).ColumnWidth = 12.833333333333334
$ws.Columns.Item(8).ColumnWidth = 11.833333333333334
